# Update "paises" (countries) data sheet + timestamp, and insert Kuwait's
# refreshed figures into its new ranking position (between Sudafrica and
# Corea del Sur), which pushes Corea del Sur / Republica Dominicana /
# Dinamarca down one row each (their own totals are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 14:05"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Straightforward per-country stat refreshes ------------------------
Set-Row 4   1408745 109 296746 1028550 16473 24 83449   # Estados Unidos
Set-Row 11  173327  156 148700 16871   1539  18 7756    # Alemania
Set-Row 15  75048   756 24900  47708   0     25 2440    # India
Set-Row 19  43211   227 0      37399   463   52 5562    # Paises Bajos
Set-Row 26  28132   219 3182   23775   103   12 1175    # Portugal
Set-Row 27  27909   637 4971   19478   351   147 3460   # Suecia
Set-Row 58  6054    51  4300   1470    35    9  284     # Finlandia
Set-Row 125 446     5   124    315     1     0  7       # Zambia

# --- Kuwait re-ranking ---------------------------------------------------
# Before: row44=Corea del Sur, row45=Republica Dominicana, row46=Dinamarca,
#         row47=Kuwait.
# After:  row44=Kuwait (new totals), row45=Corea del Sur, row46=Republica
#         Dominicana, row47=Dinamarca. Corea del Sur / Republica Dominicana
#         / Dinamarca keep their own existing totals -- only the row they
#         sit on shifts down by one.

$corea     = $ws.Range("B44:H44").Value2
$repdom    = $ws.Range("B45:H45").Value2
$dinamarca = $ws.Range("B46:H46").Value2

$ws.Cells.Item(44, 1).Value = "Kuwait"
Set-Row 44 11028 751 3263 7683 169 7 82

$ws.Cells.Item(45, 1).Value = "Corea del Sur"
$ws.Range("B45:H45").Value = $corea

$ws.Cells.Item(46, 1).Value = "Republica Dominicana"
$ws.Range("B46:H46").Value = $repdom

$ws.Cells.Item(47, 1).Value = "Dinamarca"
$ws.Range("B47:H47").Value = $dinamarca
